$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B7").Value = 1993158.79
$ws.Range("C7").Value = -55.8785068537943
$ws.Range("D7").Value = 1906
$ws.Range("E7").Value = 1906
$ws.Range("F7").Value = 1045.728641133263
$ws.Range("G7").Value = 7.91941293159022
